$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.895.83'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.667.55'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.42'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.523'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0622'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0893'
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('D12').Value = '1.904.47'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '1.648.72'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.08'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '26.900.91'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '234.45'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.00'
$ws.Range('E19').Value = '  +1.67%  '
$ws.Range('D20').Value = '0.0₃0731'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.40'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.13'
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.67'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.13'
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.87'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0494'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').Value = '1.449.19'
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.902'
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.30'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.987'
$ws.Range('E43').Value = '  +7.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.99'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('D45').Value = '1.813.65'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.782'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('D49').Value = '0.0₆0103'
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('E50').Value = '  +3.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0507'
$ws.Range('E51').Value = '  -0.01%  '
